$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT (cell type "str"/text,
# not coerced to a Number) without leaving a formula behind and without
# bumping the cell's style index. We do this by writing a string-literal
# formula, then flattening it back to a plain value with PasteSpecial
# (xlPasteValues = -4163).
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Qty executed upto date (column C) — plain numeric cells
$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 21
$ws.Range("C10").Value = 12
$ws.Range("C11").Value = 69
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 73
$ws.Range("C14").Value = 58
$ws.Range("C15").Value = 72
$ws.Range("C16").Value = 67
$ws.Range("C17").Value = 83

# Upto date Amount (column G) / Amount Since prev bill (column H) —
# stored as text (e.g. "5376.00"), so use the text-preserving helper.
Set-TextValue $ws.Range("G9")  "5376.00"
Set-TextValue $ws.Range("G10") "5664.00"
Set-TextValue $ws.Range("G11") "45678.00"
Set-TextValue $ws.Range("G13") "9928.00"
Set-TextValue $ws.Range("G14") "1334.00"

Set-TextValue $ws.Range("G19") "67980.00"
Set-TextValue $ws.Range("H19") "67980.00"
Set-TextValue $ws.Range("G21") "67980.00"
Set-TextValue $ws.Range("H21") "67980.00"
